# Updated symbol list on Fri Jan  6 06:19:48 UTC 2023 with GitHub Actions
# Refresh the Price (D), Volume(1h) (E) and Hora (G) columns on Sheet1 with
# the latest crypto snapshot. All three columns are stored as plain text in
# this sheet, so each new value is written with a leading apostrophe
# (quote-prefix) to force Excel to keep it as text instead of silently
# re-interpreting it as a number / percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{Cell="D2"; Value="256.46"},
    @{Cell="E2"; Value="-0.69%"},
    @{Cell="G2"; Value="6"},
    @{Cell="D3"; Value="27.01"},
    @{Cell="E3"; Value="0.66%"},
    @{Cell="G3"; Value="6"},
    @{Cell="D4"; Value="4.376"},
    @{Cell="E4"; Value="-7.34%"},
    @{Cell="G4"; Value="6"},
    @{Cell="D5"; Value="0.05887"},
    @{Cell="E5"; Value="-1.42%"},
    @{Cell="G5"; Value="6"},
    @{Cell="D6"; Value="6.634"},
    @{Cell="E6"; Value="-0.76%"},
    @{Cell="G6"; Value="6"},
    @{Cell="D7"; Value="0.8521"},
    @{Cell="E7"; Value="-2.43%"},
    @{Cell="G7"; Value="6"},
    @{Cell="D8"; Value="0.9389"},
    @{Cell="E8"; Value="-1.77%"},
    @{Cell="G8"; Value="6"},
    @{Cell="D9"; Value="0.01042"},
    @{Cell="E9"; Value="1,603.05%"},
    @{Cell="G9"; Value="6"},
    @{Cell="E10"; Value="-2.05%"},
    @{Cell="G10"; Value="6"},
    @{Cell="D11"; Value="0.04925"},
    @{Cell="E11"; Value="36.50%"},
    @{Cell="G11"; Value="6"},
    @{Cell="D12"; Value="0.07081"},
    @{Cell="E12"; Value="-1.50%"},
    @{Cell="G12"; Value="6"},
    @{Cell="D13"; Value="0.03075"},
    @{Cell="E13"; Value="-2.07%"},
    @{Cell="G13"; Value="6"},
    @{Cell="D14"; Value="0.09108"},
    @{Cell="E14"; Value="-1.40%"},
    @{Cell="G14"; Value="6"},
    @{Cell="D15"; Value="0.001537"},
    @{Cell="E15"; Value="-0.61%"},
    @{Cell="G15"; Value="6"},
    @{Cell="D16"; Value="0.006170"},
    @{Cell="E16"; Value="3.12%"},
    @{Cell="G16"; Value="6"},
    @{Cell="D17"; Value="3.485"},
    @{Cell="E17"; Value="-0.06%"},
    @{Cell="G17"; Value="6"},
    @{Cell="D18"; Value="3.169"},
    @{Cell="E18"; Value="-1.14%"},
    @{Cell="G18"; Value="6"},
    @{Cell="D19"; Value="2.204"},
    @{Cell="E19"; Value="-0.68%"},
    @{Cell="G19"; Value="6"},
    @{Cell="D20"; Value="0.3051"},
    @{Cell="E20"; Value="-2.04%"},
    @{Cell="G20"; Value="6"},
    @{Cell="D21"; Value="0.1269"},
    @{Cell="E21"; Value="-1.61%"},
    @{Cell="G21"; Value="6"},
    @{Cell="D22"; Value="3.917"},
    @{Cell="E22"; Value="11.16%"},
    @{Cell="G22"; Value="6"},
    @{Cell="D23"; Value="0.04264"},
    @{Cell="E23"; Value="1.02%"},
    @{Cell="G23"; Value="6"},
    @{Cell="D24"; Value="0.001221"},
    @{Cell="E24"; Value="0.08%"},
    @{Cell="G24"; Value="6"},
    @{Cell="D25"; Value="0.004281"},
    @{Cell="E25"; Value="-5.12%"},
    @{Cell="G25"; Value="6"},
    @{Cell="E26"; Value="0.00%"},
    @{Cell="G26"; Value="6"},
    @{Cell="D27"; Value="0.0001523"},
    @{Cell="E27"; Value="2.04%"},
    @{Cell="G27"; Value="6"},
    @{Cell="G28"; Value="6"},
    @{Cell="G29"; Value="6"},
    @{Cell="G30"; Value="6"},
    @{Cell="G31"; Value="6"},
    @{Cell="G32"; Value="6"},
    @{Cell="G33"; Value="6"},
    @{Cell="G34"; Value="6"},
    @{Cell="G35"; Value="6"},
    @{Cell="G36"; Value="6"},
    @{Cell="G37"; Value="6"},
    @{Cell="G38"; Value="6"},
    @{Cell="G39"; Value="6"},
    @{Cell="D40"; Value="0.03814"},
    @{Cell="E40"; Value="-0.60%"},
    @{Cell="G40"; Value="6"},
    @{Cell="D41"; Value="0.006240"},
    @{Cell="E41"; Value="4.35%"},
    @{Cell="G41"; Value="6"},
    @{Cell="D42"; Value="0.1099"},
    @{Cell="E42"; Value="-0.41%"},
    @{Cell="G42"; Value="6"},
    @{Cell="E43"; Value="0.00%"},
    @{Cell="G43"; Value="6"},
    @{Cell="E44"; Value="29.55%"},
    @{Cell="G44"; Value="6"},
    @{Cell="D45"; Value="0.00005377"},
    @{Cell="E45"; Value="-2.13%"},
    @{Cell="G45"; Value="6"},
    @{Cell="E46"; Value="0.00%"},
    @{Cell="G46"; Value="6"},
    @{Cell="D47"; Value="0.05699"},
    @{Cell="E47"; Value="-33.33%"},
    @{Cell="G47"; Value="6"},
    @{Cell="E48"; Value="11,716.51%"},
    @{Cell="G48"; Value="6"},
    @{Cell="E49"; Value="0.00%"},
    @{Cell="G49"; Value="6"},
    @{Cell="E50"; Value="0.00%"},
    @{Cell="G50"; Value="6"},
    @{Cell="G51"; Value="6"}
)

foreach ($update in $cellUpdates) {
    $ws.Range($update.Cell).Value = "'" + $update.Value
}
